# Update the TPM-derived metrics for rows 2-10 (columns G..T) on the active
# sheet, reflecting the new TPM values used to recompute the LR-pair stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => values for columns G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$newValues = @{
    2  = @(24.999262, 74.997786, 0.2094245171924971, 0.209424517192497, 3, 1, 0.2570643333333333, 0.771193, 0.1391871519274401, 0.1391871519274401, 6.426418619855334, 57.83776757869801, 0.02914920209180287, 0.02914920209180287)
    3  = @(24.999262, 74.997786, 0.2094245171924971, 0.209424517192497, 3, 1, 0.9912810000000002, 2.973843, 0.5367278196889161, 0.536727819688916, 24.78129343462201, 223.0316409115981, 0.1124039645021329, 0.1124039645021328)
    4  = @(24.999262, 74.997786, 0.2094245171924971, 0.209424517192497, 3, 1, 0.5985516666666667, 1.795655, 0.3240850283836438, 0.3240850283836438, 14.96334993553667, 134.67014941983, 0.06787135059856131, 0.06787135059856129)
    5  = @(62.40815866666667, 187.224476, 0.522807373179233, 0.5228073731792329, 3, 1, 0.2570643333333333, 0.771193, 0.1391871519274401, 0.1391871519274401, 16.04291170220756, 144.386205319868, 0.07276806927948376, 0.07276806927948375)
    6  = @(62.40815866666667, 187.224476, 0.522807373179233, 0.5228073731792329, 3, 1, 0.9912810000000002, 2.973843, 0.5367278196889161, 0.536727819688916, 61.86402193125202, 556.7761973812682, 0.2806052615237792, 0.2806052615237791)
    7  = @(62.40815866666667, 187.224476, 0.522807373179233, 0.5228073731792329, 3, 1, 0.5985516666666667, 1.795655, 0.3240850283836438, 0.3240850283836438, 37.35450738353111, 336.19056645178, 0.16943404237597, 0.1694340423759699)
    8  = @(31.96380833333333, 95.891425, 0.2677681096282701, 0.2677681096282701, 3, 1, 0.2570643333333333, 0.771193, 0.1391871519274401, 0.1391871519274401, 8.216755080002779, 73.950795720025, 0.03726988055615345, 0.03726988055615345)
    9  = @(31.96380833333333, 95.891425, 0.2677681096282701, 0.2677681096282701, 3, 1, 0.9912810000000002, 2.973843, 0.5367278196889161, 0.536727819688916, 31.68511588847501, 285.1660429962751, 0.1437185936630041, 0.143718593663004)
    10 = @(31.96380833333333, 95.891425, 0.2677681096282701, 0.2677681096282701, 3, 1, 0.5985516666666667, 1.795655, 0.3240850283836438, 0.3240850283836438, 19.13199075093056, 172.187916758375, 0.08677963540911256, 0.08677963540911254)
}

$columns = @("G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $ws.Range("$col$row").Value = $vals[$i]
    }
}
